$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") '61.525.87'
Set-TextCell $ws.Range("E2") '  +0.69%  '

# Row 3
Set-TextCell $ws.Range("D3") '3.446.09'
Set-TextCell $ws.Range("E3") '  +1.54%  '

# Row 4
Set-TextCell $ws.Range("D4") '0.999'
Set-TextCell $ws.Range("E4") '  -0.08%  '

# Row 5
Set-TextCell $ws.Range("D5") '577.35'
Set-TextCell $ws.Range("E5") '  +0.84%  '

# Row 6
Set-TextCell $ws.Range("D6") '145.34'
Set-TextCell $ws.Range("E6") '  +4.98%  '

# Row 7
Set-TextCell $ws.Range("D7") '3.446.40'
Set-TextCell $ws.Range("E7") '  +1.56%  '

# Row 8
Set-TextCell $ws.Range("E8") '  +0.03%  '

# Row 9
Set-TextCell $ws.Range("E9") '  +2.52%  '

# Row 10
Set-TextCell $ws.Range("D10") '7.70'
Set-TextCell $ws.Range("E10") '  -0.04%  '

# Row 11
Set-TextCell $ws.Range("E11") '  +3.86%  '

# Row 12
Set-TextCell $ws.Range("D12") '0.390'
Set-TextCell $ws.Range("E12") '  +2.81%  '

# Row 13
Set-TextCell $ws.Range("D13") '4.035.12'
Set-TextCell $ws.Range("E13") '  +1.56%  '

# Row 14
Set-TextCell $ws.Range("D14") '28.57'
Set-TextCell $ws.Range("E14") '  +7.29%  '

# Row 15
Set-TextCell $ws.Range("E15") '  -0.38%  '

# Row 16
Set-TextCell $ws.Range("E16") '  +1.44%  '

# Row 17
Set-TextCell $ws.Range("D17") '3.446.52'
Set-TextCell $ws.Range("E17") '  +1.45%  '

# Row 18
Set-TextCell $ws.Range("D18") '61.647.44'
Set-TextCell $ws.Range("E18") '  +0.80%  '

# Row 19
Set-TextCell $ws.Range("D19") '6.35'
Set-TextCell $ws.Range("E19") '  +6.94%  '

# Row 20
Set-TextCell $ws.Range("E20") '  +3.88%  '

# Row 21
Set-TextCell $ws.Range("E21") '  +1.45%  '

# Row 22
Set-TextCell $ws.Range("D22") '403.10'
Set-TextCell $ws.Range("E22") '  +7.71%  '

# Row 23
Set-TextCell $ws.Range("E23") '  +3.36%  '

# Row 24
Set-TextCell $ws.Range("D24") '74.49'
Set-TextCell $ws.Range("E24") '  +4.80%  '

# Row 25
Set-TextCell $ws.Range("B25") 'LEO'
Set-TextCell $ws.Range("C25") 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell $ws.Range("D25") '5.77'
Set-TextCell $ws.Range("E25") '  +0.68%  '

# Row 26
Set-TextCell $ws.Range("B26") 'Dai'
Set-TextCell $ws.Range("C26") 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws.Range("D26") '0.998'
Set-TextCell $ws.Range("E26") '  -0.26%  '

# Row 27
Set-TextCell $ws.Range("E27") '  +1.10%  '

# Row 28
Set-TextCell $ws.Range("D28") '3.588.57'
Set-TextCell $ws.Range("E28") '  +1.58%  '

# Row 29
Set-TextCell $ws.Range("E29") '  +4.88%  '

# Row 30
Set-TextCell $ws.Range("D30") '7.65'
Set-TextCell $ws.Range("E30") '  +3.82%  '

# Row 31
Set-TextCell $ws.Range("E31") '  +0.05%  '

# Row 32
Set-TextCell $ws.Range("E32") '  +2.39%  '

# Row 33
Set-TextCell $ws.Range("E33") '  +2.23%  '

# Row 34
Set-TextCell $ws.Range("D34") '1.45'

# Row 36
Set-TextCell $ws.Range("D36") '23.92'
Set-TextCell $ws.Range("E36") '  +2.32%  '

# Row 37
Set-TextCell $ws.Range("B37") 'Aptos'
Set-TextCell $ws.Range("C37") 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws.Range("D37") '7.05'
Set-TextCell $ws.Range("E37") '  +3.19%  '

# Row 38
Set-TextCell $ws.Range("B38") 'RenzoRestakedETH'
Set-TextCell $ws.Range("C38") 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextCell $ws.Range("D38") '3.473.83'
Set-TextCell $ws.Range("E38") '  +1.77%  '

# Row 39
Set-TextCell $ws.Range("B39") 'NEARProtocol'
Set-TextCell $ws.Range("C39") 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws.Range("D39") '5.15'
Set-TextCell $ws.Range("E39") '  +1.24%  '

# Row 40
Set-TextCell $ws.Range("B40") 'ImmutableX'
Set-TextCell $ws.Range("C40") 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws.Range("D40") '1.56'
Set-TextCell $ws.Range("E40") '  +0.40%  '

# Row 41
Set-TextCell $ws.Range("D41") '166.99'
Set-TextCell $ws.Range("E41") '  +0.29%  '

# Row 42
Set-TextCell $ws.Range("D42") '0.0792'
Set-TextCell $ws.Range("E42") '  +3.27%  '

# Row 43
Set-TextCell $ws.Range("D43") '27.23'
Set-TextCell $ws.Range("E43") '  +5.09%  '

# Row 44
Set-TextCell $ws.Range("E44") '  +3.44%  '

# Row 45
Set-TextCell $ws.Range("D45") '4.53'
Set-TextCell $ws.Range("E45") '  +3.37%  '

# Row 46
Set-TextCell $ws.Range("B46") 'Stacks'
Set-TextCell $ws.Range("C46") 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell $ws.Range("D46") '1.74'
Set-TextCell $ws.Range("E46") '  -0.10%  '

# Row 47
Set-TextCell $ws.Range("B47") 'OKB'
Set-TextCell $ws.Range("C47") 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell $ws.Range("D47") '42.45'
Set-TextCell $ws.Range("E47") '  +1.41%  '

# Row 48
Set-TextCell $ws.Range("B48") 'FirstDigitalUSD'
Set-TextCell $ws.Range("C48") 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell $ws.Range("D48") '1.00'
Set-TextCell $ws.Range("E48") '  -0.08%  '

# Row 49
Set-TextCell $ws.Range("D49") '2.611.48'
Set-TextCell $ws.Range("E49") '  +4.06%  '

# Row 50
Set-TextCell $ws.Range("E50") '  -1.40%  '

# Row 51
Set-TextCell $ws.Range("E51") '  +2.63%  '
